$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.134.00'
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').Value = '1.782.40'
$ws.Range('E3').Value = '  -1.94%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '337.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3857'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3425'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.86'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.188'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07446'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.62'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.426'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.59%  '
$ws.Range('D15').Value = '1.782.82'
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.110'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001090'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06654'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.30'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.57%  '
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.48'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.513'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').Value = '27.125.70'
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('E24').Value = '  -6.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.364'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.02%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.51%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.492'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.438'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '155.74'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('D30').Value = '1.985.59'
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '133.85'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.978'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.988'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08680'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.98'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.82%  '
$ws.Range('E36').Value = '  -3.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.385'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6813'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06318'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02338'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2181'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.71%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.237'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.07%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.412'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.27'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6393'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.856'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.188'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '131.36'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07109'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.09'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.80%  '
